$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A1 holds the sheet/title label (shared string "HK_G_acc_G") - rewritten as part of the refresh
$ws.Range("A1").Value = "HK_G_acc_G"

# Updated accuracy values (A2:A49) from the latest threshold run
$ws.Range("A2").Value = 52.072072072072075
$ws.Range("A3").Value = 51.711711711711715
$ws.Range("A4").Value = 51.891891891891895
$ws.Range("A5").Value = 51.531531531531527
$ws.Range("A6").Value = 51.351351351351347
$ws.Range("A7").Value = 51.531531531531527
$ws.Range("A8").Value = 53.693693693693689
$ws.Range("A9").Value = 53.153153153153156
$ws.Range("A10").Value = 53.153153153153156
$ws.Range("A11").Value = 53.333333333333336
$ws.Range("A12").Value = 52.432432432432428
$ws.Range("A13").Value = 52.792792792792795
$ws.Range("A14").Value = 54.054054054054056
$ws.Range("A15").Value = 54.054054054054056
$ws.Range("A16").Value = 53.873873873873876
$ws.Range("A17").Value = 52.612612612612608
$ws.Range("A18").Value = 52.792792792792795
$ws.Range("A19").Value = 52.252252252252248
$ws.Range("A20").Value = 52.612612612612608
$ws.Range("A21").Value = 52.612612612612608
$ws.Range("A22").Value = 52.972972972972975
$ws.Range("A23").Value = 50.810810810810814
$ws.Range("A24").Value = 49.909909909909913
$ws.Range("A25").Value = 50.990990990990994
$ws.Range("A26").Value = 52.792792792792795
$ws.Range("A27").Value = 52.612612612612608
$ws.Range("A28").Value = 53.513513513513509
$ws.Range("A29").Value = 54.414414414414416
$ws.Range("A30").Value = 53.873873873873876
$ws.Range("A31").Value = 54.234234234234236
$ws.Range("A32").Value = 50.450450450450447
$ws.Range("A33").Value = 51.531531531531527
$ws.Range("A34").Value = 51.711711711711715
$ws.Range("A35").Value = 52.072072072072075
$ws.Range("A36").Value = 51.531531531531527
$ws.Range("A37").Value = 56.216216216216218
$ws.Range("A38").Value = 50.810810810810814
$ws.Range("A39").Value = 51.711711711711715
$ws.Range("A40").Value = 51.171171171171167
$ws.Range("A41").Value = 52.612612612612608
$ws.Range("A42").Value = 53.153153153153156
$ws.Range("A43").Value = 53.873873873873876
$ws.Range("A44").Value = 53.333333333333336
$ws.Range("A45").Value = 52.252252252252248
$ws.Range("A46").Value = 52.432432432432428
$ws.Range("A47").Value = 51.531531531531527
$ws.Range("A48").Value = 54.054054054054056
$ws.Range("A49").Value = 52.972972972972975
